# SMPTE test pattern workbook update
# Adds support for SMPTE 16x9, SMPTE 4x3 and EBU test patterns and moves the
# extra EBU colors into the shared-strings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# ---------------------------------------------------------------------------
# SMPTE 4x3 block (rows 16-19 use B14/C14 as their reference row instead of
# the B2/C2 reference used by rows 4-13)
# ---------------------------------------------------------------------------
$ws.Range("D16:D19").Formula = '=B16/$B$14'
$ws.Range("E16:E19").Formula = '=C16/$C$14'

# ---------------------------------------------------------------------------
# Row 20 no longer carries a formula
# ---------------------------------------------------------------------------
$ws.Range("D20:E20").ClearContents()

# ---------------------------------------------------------------------------
# New EBU test pattern reference row (21) and first data row (22)
# ---------------------------------------------------------------------------
$ws.Range("D21:E21").ClearContents()
$ws.Range("B21").Value2 = 150
$ws.Range("C21").Value2 = 200

$ws.Range("B22").Value2 = 150
$ws.Range("C22").Value2 = 25
$ws.Range("D22").Formula = '=B22/$B$21'
$ws.Range("E22").Formula = '=C22/$C$21'

# ---------------------------------------------------------------------------
# EBU block (rows 23-26) use B21/C21 as their reference
# ---------------------------------------------------------------------------
$ws.Range("D23:D26").Formula = '=B23/$B$21'
$ws.Range("E23:E26").Formula = '=C23/$C$21'

# ---------------------------------------------------------------------------
# New EBU color strings, referenced from column I (rows 25-30)
# ---------------------------------------------------------------------------
$ws.Range("I25").Value2 = "BFBFBF"
$ws.Range("I26").Value2 = "BFBF00"
$ws.Range("I27").Value2 = "00BFBF"
$ws.Range("I28").Value2 = "00BF00"
$ws.Range("I29").Value2 = "BF0000"
$ws.Range("I30").Value2 = "0000BF"
$ws.Range("I31").Value2 = 0

# ---------------------------------------------------------------------------
# Move the active selection
# ---------------------------------------------------------------------------
[void]$ws.Range("I32").Select()
